$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. site_data: area_plan (C24) changes from 1 to 0.096.
# ---------------------------------------------------------------------
$wsSite = $wb.Worksheets.Item("site_data")
$wsSite.Range("C24").Value = 0.096

# ---------------------------------------------------------------------
# 2. Preserve the full, original "surface_data" content on a new sheet
#    named "surface_data_2" (moved to the end of the workbook).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("surface_data").Name = "surface_data_2"
$wb.Worksheets.Item("surface_data_2").Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# ---------------------------------------------------------------------
# 3. Re-create "surface_data" as a short summary sheet (in the same
#    slot, right before "pollution_data").
# ---------------------------------------------------------------------
$wb.Worksheets.Item("pollution_data").Select()
$wsSurface = $wb.Worksheets.Add()
$wsSurface.Name = "surface_data"

$wsSurfaceFull = $wb.Worksheets.Item("surface_data_2")
$wsSurfaceFull.Range("A1:E4").Copy($wsSurface.Range("A1"))

# Copy() drops the formula text (only the cached value survives) -
# restore it explicitly.
$wsSurface.Range("C2").Formula = "=area_plan*100"

# Match the workbook's usual top/bottom page margins (0.7/0.7 left-
# right and 0.3/0.3 header/footer are already the engine defaults).
$wsSurface.PageSetup.TopMargin = 56.692913399999995
$wsSurface.PageSetup.BottomMargin = 56.692913399999995

# New summary values.
$wsSurface.Range("C3").Value = 9.6
$wsSurface.Range("D4").Value = 100
$wsSurface.Range("E4").Value = 0.9

# The rectangular Copy() above stamped out explicit-but-blank cells
# (B3/D3/E3, C4) - drop them so empty cells stay implicit again.
$wsSurface.Range("B3").ClearContents()
$wsSurface.Range("D3").ClearContents()
$wsSurface.Range("E3").ClearContents()
$wsSurface.Range("C4").ClearContents()

# ---------------------------------------------------------------------
# 4. Selections / active sheet to mirror the saved UI state.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("pollution_data").Range("G22").Select()
$wsSurfaceFull.Range("A1:E4").Select()
$wsSite.Range("C24").Select()

$wsSurface.Range("D5").Select()
$wsSurface.Activate()

$wb.RecalcAll()
